# Tuxedo Round 1 results
# Fill in the Round 1 winner column (O) with the same team names already
# present in the adjacent "advance" column (P), for each of the 8 first
# round matchups (rows 2,4,6,8,10,12,14,16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(2, 4, 6, 8, 10, 12, 14, 16)

foreach ($r in $rows) {
    $winner = $ws.Range("P$r").Value2
    $ws.Range("O$r").Value = $winner
}

# Move the active selection to O16, matching the saved cursor position.
$ws.Range("O16").Select()
